# Apply volatility_comparison value updates to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: Future_Price ---
$ws.Range("B5").Value = 331.89844
$ws.Range("C5").Value = 110.63281
$ws.Range("D5").Value = 110.63281
$ws.Range("E5").Value = 110.63281
$ws.Range("F5").Value = 110.6171875
$ws.Range("G5").Value = 110.6171875
$ws.Range("H5").Value = 110.6171875
$ws.Range("I5").Value = 110.6171875
$ws.Range("J5").Value = 110.6171875
$ws.Range("K5").Value = 110.6171875

# --- Row 6: Days ---
$ws.Range("B6").Value = 8.406249999999993
$ws.Range("C6").Value = 0.1354199999999865
$ws.Range("D6").Value = 3.135419999999994
$ws.Range("E6").Value = 5.135419999999987
$ws.Range("F6").Value = 0.1
$ws.Range("G6").Value = 0.1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3

# --- Row 7: Time_Years ---
$ws.Range("B7").Value = 0.0230308219178082
$ws.Range("C7").Value = 0.0003710136986301
$ws.Range("D7").Value = 0.008590191780821901
$ws.Range("E7").Value = 0.0140696438356164
$ws.Range("F7").Value = 0.0003968253968253968
$ws.Range("G7").Value = 0.0003968253968253968
$ws.Range("H7").Value = 0.003968253968253968
$ws.Range("I7").Value = 0.003968253968253968
$ws.Range("J7").Value = 0.0119047619047619
$ws.Range("K7").Value = 0.0119047619047619

# --- Row 8: Market_Price ---
$ws.Range("F8").Value = 1.2
$ws.Range("G8").Value = 8.5
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 17
$ws.Range("J8").Value = 17
$ws.Range("K8").Value = 24.5

# --- Row 9: Market_Price_64ths (text values - force text storage so the
#     numeric-looking strings aren't auto-coerced into numbers) ---
$row9 = $ws.Range("F9:K9")
$row9Format = $row9.NumberFormat
$row9.NumberFormat = "@"
$ws.Range("F9").Value = "1.20"
$ws.Range("G9").Value = "8.50"
$ws.Range("H9").Value = "10.00"
$ws.Range("I9").Value = "17.00"
$ws.Range("J9").Value = "17.00"
$ws.Range("K9").Value = "24.50"
$row9.NumberFormat = $row9Format

# --- Row 10: Bid ---
$ws.Range("G10").Value = 8
$ws.Range("H10").Value = 9
$ws.Range("I10").Value = 16
$ws.Range("J10").Value = 16
$ws.Range("K10").Value = 24

# --- Row 11: Ask ---
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 11
$ws.Range("I11").Value = 18
$ws.Range("J11").Value = 18
$ws.Range("K11").Value = 25

# --- Row 12: Market_Vol ---
$ws.Range("B12").Value = 23.70031
$ws.Range("D12").Value = 7.2166
$ws.Range("E12").Value = 7.72159
$ws.Range("F12").Value = 6.59
$ws.Range("G12").Value = 6.68
$ws.Range("H12").Value = 8.59
$ws.Range("K12").Value = 7.36

# --- Row 13: Calculated_Vol (text values - same text-forcing trick) ---
$row13 = $ws.Range("F13:K13")
$row13Format = $row13.NumberFormat
$row13.NumberFormat = "@"
$ws.Range("F13").Value = "8.12"
$ws.Range("G13").Value = "7.00"
$ws.Range("H13").Value = "8.60"
$ws.Range("I13").Value = "8.02"
$ws.Range("J13").Value = "7.53"
$ws.Range("K13").Value = "7.37"
$row13.NumberFormat = $row13Format

# --- Row 14: Difference_from_Calculated ---
$ws.Range("F14").Value = -1.529999999999999
$ws.Range("G14").Value = -0.3200000000000003
$ws.Range("H14").Value = -0.009999999999999787
$ws.Range("I14").Value = -0.9699999999999998
$ws.Range("J14").Value = -1.06
$ws.Range("K14").Value = -0.009999999999999787
